# Auto-generated edit script: updates cryptos.xlsx price/volume data
# per the commit "Updated cryptos list on Thu Mar 16 10:27:32 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '24.807.50'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.658.23'
$ws.Cells.Item(3, 5).Value = '  -2.39%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.000'
$ws.Cells.Item(4, 5).Value = '  -0.78%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''316.76'
$ws.Cells.Item(5, 5).Value = '  +1.72%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''0.9975'
$ws.Cells.Item(6, 5).Value = '  -0.48%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.3630'
$ws.Cells.Item(7, 5).Value = '  -2.45%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''46.76'
$ws.Cells.Item(8, 5).Value = '  -4.79%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.3265'
$ws.Cells.Item(9, 5).Value = '  -4.23%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -5.55%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.07059'
$ws.Cells.Item(11, 5).Value = '  -5.19%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.9974'
$ws.Cells.Item(12, 5).Value = '  -0.59%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''6.044'
$ws.Cells.Item(13, 5).Value = '  -4.29%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -6.08%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '1.664.46'
$ws.Cells.Item(15, 5).Value = '  -2.01%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''6.621'
$ws.Cells.Item(16, 5).Value = '  -5.13%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''0.00001048'
$ws.Cells.Item(17, 5).Value = '  -6.24%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.06621'
$ws.Cells.Item(18, 5).Value = '  -1.15%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''0.9974'
$ws.Cells.Item(19, 5).Value = '  -0.44%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''79.16'
$ws.Cells.Item(20, 5).Value = '  -4.83%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''5.923'
$ws.Cells.Item(21, 5).Value = '  -6.12%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''15.75'
$ws.Cells.Item(22, 5).Value = '  -7.96%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''12.55'
$ws.Cells.Item(23, 5).Value = '  -2.66%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '24.776.95'
$ws.Cells.Item(24, 5).Value = '  +0.15%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''2.432'
$ws.Cells.Item(25, 5).Value = '  -1.10%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''2.401'
$ws.Cells.Item(26, 5).Value = '  -12.85%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''148.67'
$ws.Cells.Item(27, 5).Value = '  -0.11%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''18.60'
$ws.Cells.Item(28, 5).Value = '  -7.72%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'ImmutableX'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(29, 4).Value = '''1.225'
$ws.Cells.Item(29, 5).Value = '  -0.73%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(30, 4).Value = '1.847.59'
$ws.Cells.Item(30, 5).Value = '  -1.96%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''125.59'
$ws.Cells.Item(31, 5).Value = '  -4.48%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''4.082'
$ws.Cells.Item(32, 5).Value = '  -3.09%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''5.827'
$ws.Cells.Item(33, 5).Value = '  -13.01%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''0.08452'
$ws.Cells.Item(34, 5).Value = '  -2.90%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''1.672'
$ws.Cells.Item(35, 5).Value = '  -4.66%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''12.27'
$ws.Cells.Item(36, 5).Value = '  -9.01%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +1.52%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''5.198'
$ws.Cells.Item(38, 5).Value = '  -5.91%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).Value = '''0.06044'
$ws.Cells.Item(39, 5).Value = '  -8.60%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.02239'
$ws.Cells.Item(40, 5).Value = '  -6.52%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.2073'
$ws.Cells.Item(41, 5).Value = '  -6.35%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''8.226'
$ws.Cells.Item(42, 5).Value = '  -8.54%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.9968'
$ws.Cells.Item(43, 5).Value = '  -0.51%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''0.5925'
$ws.Cells.Item(44, 5).Value = '  -7.06%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''3.818'
$ws.Cells.Item(45, 5).Value = '  +0.07%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''12.82'
$ws.Cells.Item(46, 5).Value = '  -6.35%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.5648'
$ws.Cells.Item(47, 5).Value = '  -7.07%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''125.04'
$ws.Cells.Item(48, 5).Value = '  -2.80%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''1.955'
$ws.Cells.Item(49, 5).Value = '  -7.12%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.06993'
$ws.Cells.Item(50, 5).Value = '  -3.44%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''1.197'
$ws.Cells.Item(51, 5).Value = '  -1.43%  '

# Normalize style on touched D-column cells (strip any quote-prefix styling
# introduced by assigning numeric-looking text so the cells keep the sheet
# default style, matching the source data's formatting).
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(51, 4).ClearFormats()
